# "sua loi export excel" - fix the Excel export template:
# the search-criteria block on the sheet was missing a "Linh vuc:" (Field /
# Sector) filter row. Add it as row 5 and push the existing criteria rows
# (previously rows 5-10: Tinh/TP:, Quan/Huyen:, Xa/Phuong:, Thoi gian tiep
# nhan:, Giai doan:, Ket qua:) down into rows 6-11, which were already
# blank, so the rest of the sheet (header/data rows 14-15) does not move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift D5:E10 down into D6:E11, bottom row first so we never overwrite a
# source cell before it has been copied from.
for ($r = 10; $r -ge 5; $r--) {
    $dst = $r + 1

    $ws.Range("D$r").Copy()
    $ws.Range("D$dst").PasteSpecial(-4122)
    $ws.Range("D$dst").Value2 = $ws.Range("D$r").Value2

    $ws.Range("E$r").Copy()
    $ws.Range("E$dst").PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# Row 5 now still holds the old "Tinh/TP:" text (copied onto row 6 above) -
# turn it into the new "Linh vuc:" label, matching the other label cells'
# formatting, and drop its (unused) companion input cell in column E.
$ws.Range("D5").Value2 = "Lĩnh vực:"
$ws.Range("E5").Clear()

# Leave the selection where a user would land after adding this row.
$ws.Range("D11").Select()
